$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering / content for the "Organization" column (A2:A31).
# rank_group values in column B are unchanged by this edit.
$orgs = @(
    "Owl Ventures",
    "Bill & Melinda Gates Foundation",
    "Y Combinator",
    "Reach Capital",
    "Google for Education",
    "U.S. Dept. of Education (EIR Program)",
    "National Science Foundation (NSF)",
    "500 Global Flagship VC (non-accelerator checks)",
    "IES SBIR (ED/IES)",
    "500 Global (seed/accelerator)",
    "Buffalo Sabres Foundation",
    "TGR Foundation (Tiger Woods)",
    "Berkeley SkyDeck Fund (UC Berkeley)",
    "Penn State University - Outreach & Engagement",
    "Chan Zuckerberg Initiative (CZI)",
    "Nashville Predators Foundation",
    "Portland Trail Blazers Foundation",
    "NFL Foundation",
    "San Jose Sharks Foundation",
    "San Antonio Spurs - Spurs Give",
    "Indiana Pacers Foundation",
    "Cleveland Cavaliers Community Foundation",
    "Austin FC - 4ATX Foundation",
    "Baltimore Ravens Foundation",
    "Toronto FC - MLSE Foundation",
    "Houston Texans Foundation",
    "Philadelphia Eagles Foundation",
    "Oakland Roots SC",
    "Jacksonville Jaguars Foundation",
    "Florida State University Research Foundation"
)

for ($i = 0; $i -lt $orgs.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $orgs[$i]
}
